$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.783.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "'3.796.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'701.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.67%  "
$ws.Range("D6").Value = "'172.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.41%  "
$ws.Range("D7").Value = "'3.795.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +6.65%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +7.89%  "
$ws.Range("D14").Value = "'36.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.51%  "
$ws.Range("D15").Value = "'4.436.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "'3.792.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "'70.827.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "'17.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +16.30%  "
$ws.Range("D22").Value = "'481.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "'84.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "'12.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").Value = "'2.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("D28").Value = "'10.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "'3.946.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'3.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.27%  "
$ws.Range("D32").Value = "'7.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "'29.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "'3.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("D40").Value = "'6.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.84%  "
$ws.Range("E41").Value = "  +12.27%  "
$ws.Range("D42").Value = "'0.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  +22.65%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D46").Value = "'162.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").Value = "'49.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'0.303"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  +2.40%  "
